$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear B2 and F2 (they become empty cells, removed from the sheet)
$ws.Range("B2").ClearContents()
$ws.Range("F2").ClearContents()

# Update G2 to TRUE (keep boolean type/style)
$ws.Range("G2").Value = $true

# Delete row 3 entirely (C3 was 0)
$ws.Range("A3:G3").EntireRow.Delete()

# Set selection to B2
$ws.Range("B2").Select()
